$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 45. This shifts the existing rows 45 and 46
# down to 46 and 47 respectively, matching the target layout described in
# the diff (old row45 -> new row46, old row46 -> new row47).
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly price record.
$ws.Cells.Item(45, 1).Value = 11
$ws.Cells.Item(45, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(45, 3).Value = "Bíobío"
$ws.Cells.Item(45, 4).Value = 44509
$ws.Cells.Item(45, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(45, 5).Value = 8
$ws.Cells.Item(45, 6).Value = 100112001
$ws.Cells.Item(45, 7).Value = "Berenjena"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 100
$ws.Cells.Item(45, 11).Value = 8000
$ws.Cells.Item(45, 12).Value = 9000
$ws.Cells.Item(45, 13).Value = 8500
$ws.Cells.Item(45, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(45, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(45, 16).Value = 142
$ws.Cells.Item(45, 17).Value = 60
$ws.Cells.Item(45, 18).Value = "Hortaliza"
